# Redid measurements for 2 extensor tests
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtTest10mm")

# Row 4 was previously empty; populate it in place (no row shift) with the
# "tendon" label + value.
$ws.Range("B4").Value = "tendon"
$ws.Range("B4").Font.Bold = $true
$ws.Range("C4").Value = 30

# Update the measured values for Test #1 (column C) and Test #2 (column D)
$ws.Range("C6").Value = 16.268999999999998
$ws.Range("D6").Value = 11.218999999999999

$ws.Range("C7").Value = 120
$ws.Range("D7").Value = 108.5

# C8 used to hold a formula (=90-54.6); it is now a plain measured value
$ws.Range("C8").Value = 33.6
$ws.Range("D8").Value = 36.700000000000003

$ws.Range("C9").Value = 30.5
$ws.Range("D9").Value = 30

$ws.Range("C10").Value = 518
$ws.Range("D10").Value = 520

$ws.Range("C13").Value = 38
$ws.Range("D13").Value = 35

# Move the active selection to D14, matching the post-edit cursor position
$ws.Range("D14").Select()

$wb.Application.Calculate()
